$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.428213333333333
$ws.Range("H2").Value = 4.28464
$ws.Range("I2").Value = 0.3767260624985217
$ws.Range("J2").Value = 0.3767260624985217
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.4445023333333333
$ws.Range("N2").Value = 1.333507
$ws.Range("O2").Value = 0.009977046255258984
$ws.Range("P2").Value = 0.009977046255258982
$ws.Range("Q2").Value = 0.6348441591644444
$ws.Range("R2").Value = 5.713597432479999
$ws.Range("S2").Value = 0.003758613351109338
$ws.Range("T2").Value = 0.003758613351109337

# Row 3
$ws.Range("G3").Value = 1.428213333333333
$ws.Range("H3").Value = 4.28464
$ws.Range("I3").Value = 0.3767260624985217
$ws.Range("J3").Value = 0.3767260624985217
$ws.Range("O3").Value = 0.9569553279219795
$ws.Range("P3").Value = 0.9569553279219793
$ws.Range("Q3").Value = 60.89151888940444
$ws.Range("R3").Value = 548.0236700046399
$ws.Range("S3").Value = 0.360510012675029
$ws.Range("T3").Value = 0.3605100126750289

# Row 4
$ws.Range("G4").Value = 1.428213333333333
$ws.Range("H4").Value = 4.28464
$ws.Range("I4").Value = 0.3767260624985217
$ws.Range("J4").Value = 0.3767260624985217
$ws.Range("M4").Value = 0.851471
$ws.Range("N4").Value = 2.554413
$ws.Range("O4").Value = 0.01911163320180161
$ws.Range("P4").Value = 0.01911163320180161
$ws.Range("Q4").Value = 1.216082235146666
$ws.Range("R4").Value = 10.94474011632
$ws.Range("S4").Value = 0.007199850324030737
$ws.Range("T4").Value = 0.007199850324030735

# Row 5
$ws.Range("G5").Value = 1.428213333333333
$ws.Range("H5").Value = 4.28464
$ws.Range("I5").Value = 0.3767260624985217
$ws.Range("J5").Value = 0.3767260624985217
$ws.Range("M5").Value = 0.6217743333333333
$ws.Range("N5").Value = 1.865323
$ws.Range("O5").Value = 0.01395599262095996
$ws.Range("P5").Value = 0.01395599262095996
$ws.Range("Q5").Value = 0.888026393191111
$ws.Range("R5").Value = 7.99223753872
$ws.Range("S5").Value = 0.005257586148352669
$ws.Range("T5").Value = 0.005257586148352669

# Row 6
$ws.Range("I6").Value = 0.01640263292535353
$ws.Range("J6").Value = 0.01640263292535352
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.4445023333333333
$ws.Range("N6").Value = 1.333507
$ws.Range("O6").Value = 0.009977046255258984
$ws.Range("P6").Value = 0.009977046255258982
$ws.Range("Q6").Value = 0.02764108126344444
$ws.Range("R6").Value = 0.248769731371
$ws.Range("S6").Value = 0.0001636498274042861
$ws.Range("T6").Value = 0.000163649827404286

# Row 7
$ws.Range("I7").Value = 0.01640263292535353
$ws.Range("J7").Value = 0.01640263292535352
$ws.Range("O7").Value = 0.9569553279219795
$ws.Range("P7").Value = 0.9569553279219793
$ws.Range("S7").Value = 0.01569658696986554
$ws.Range("T7").Value = 0.01569658696986554

# Row 8
$ws.Range("I8").Value = 0.01640263292535353
$ws.Range("J8").Value = 0.01640263292535352
$ws.Range("M8").Value = 0.851471
$ws.Range("N8").Value = 2.554413
$ws.Range("O8").Value = 0.01911163320180161
$ws.Range("P8").Value = 0.01911163320180161
$ws.Range("Q8").Value = 0.05294815648766667
$ws.Range("R8").Value = 0.4765334083889999
$ws.Range("S8").Value = 0.0003134811040131508
$ws.Range("T8").Value = 0.0003134811040131506

# Row 9
$ws.Range("I9").Value = 0.01640263292535353
$ws.Range("J9").Value = 0.01640263292535352
$ws.Range("M9").Value = 0.6217743333333333
$ws.Range("N9").Value = 1.865323
$ws.Range("O9").Value = 0.01395599262095996
$ws.Range("P9").Value = 0.01395599262095996
$ws.Range("Q9").Value = 0.03866462240211111
$ws.Range("R9").Value = 0.347981601619
$ws.Range("S9").Value = 0.0002289150240705487
$ws.Range("T9").Value = 0.0002289150240705486

# Row 10
$ws.Range("G10").Value = 2.300721333333334
$ws.Range("H10").Value = 6.902164000000001
$ws.Range("I10").Value = 0.6068713045761248
$ws.Range("J10").Value = 0.6068713045761248
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.4445023333333333
$ws.Range("N10").Value = 1.333507
$ws.Range("O10").Value = 0.009977046255258984
$ws.Range("P10").Value = 0.009977046255258982
$ws.Range("Q10").Value = 1.022676001016445
$ws.Range("R10").Value = 9.204084009148001
$ws.Range("S10").Value = 0.00605478307674536
$ws.Range("T10").Value = 0.006054783076745359

# Row 11
$ws.Range("G11").Value = 2.300721333333334
$ws.Range("H11").Value = 6.902164000000001
$ws.Range("I11").Value = 0.6068713045761248
$ws.Range("J11").Value = 0.6068713045761248
$ws.Range("O11").Value = 0.9569553279219795
$ws.Range("P11").Value = 0.9569553279219793
$ws.Range("Q11").Value = 98.09067963324047
$ws.Range("R11").Value = 882.8161166991641
$ws.Range("S11").Value = 0.580748728277085
$ws.Range("T11").Value = 0.5807487282770849

# Row 12
$ws.Range("G12").Value = 2.300721333333334
$ws.Range("H12").Value = 6.902164000000001
$ws.Range("I12").Value = 0.6068713045761248
$ws.Range("J12").Value = 0.6068713045761248
$ws.Range("M12").Value = 0.851471
$ws.Range("N12").Value = 2.554413
$ws.Range("O12").Value = 0.01911163320180161
$ws.Range("P12").Value = 0.01911163320180161
$ws.Range("Q12").Value = 1.958997494414667
$ws.Range("R12").Value = 17.630977449732
$ws.Range("S12").Value = 0.01159830177375773
$ws.Range("T12").Value = 0.01159830177375772

# Row 13
$ws.Range("G13").Value = 2.300721333333334
$ws.Range("H13").Value = 6.902164000000001
$ws.Range("I13").Value = 0.6068713045761248
$ws.Range("J13").Value = 0.6068713045761248
$ws.Range("M13").Value = 0.6217743333333333
$ws.Range("N13").Value = 1.865323
$ws.Range("O13").Value = 0.01395599262095996
$ws.Range("P13").Value = 0.01395599262095996
$ws.Range("Q13").Value = 1.430529473219111
$ws.Range("R13").Value = 12.874765258972
$ws.Range("S13").Value = 0.008469491448536742
$ws.Range("T13").Value = 0.008469491448536742
